# Generate Report for Handback
# - Update the "Status" text (Overview + per-language sheets) from
#   "Ready for handoff" to "Handed back: in sync with en-US"
# - Record the handback completion datetime for zh-cn and de-de
# - Add "Latest Target File" / "Latest Handback File" hyperlink cells
#   (columns F/G) for both rows of the zh-cn and de-de detail sheets

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: Status column + new Target/Handback file links ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("H2").Value = "2016-03-17 16:11:36"
$zhcn.Range("H3").Value = "2016-03-17 16:11:36"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/c2d708bfa23c3165cf7fa1ec8868c123901e405f/e2e/a1897a40-c0ff-4968-99fa-885e084b0050.md", "", "", "a1897a40-c0ff-4968-99fa-885e084b0050.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f9ca0064a146a301c80a7850a22de47c50b6f04/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/a1897a40-c0ff-4968-99fa-885e084b0050.0ffc24064d9e76e357c05896c12c1906504eaec3.zh-cn.xlf", "", "", "a1897a40-c0ff-4968-99fa-885e084b0050.0ffc24064d9e76e357c05896c12c1906504eaec3.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/c2d708bfa23c3165cf7fa1ec8868c123901e405f/e2e/e6dee80d-419c-471d-8ae0-9bd60d70567b.md", "", "", "e6dee80d-419c-471d-8ae0-9bd60d70567b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f9ca0064a146a301c80a7850a22de47c50b6f04/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/e6dee80d-419c-471d-8ae0-9bd60d70567b.60b44875e55860e846c7f9475838e5c52eca4020.zh-cn.xlf", "", "", "e6dee80d-419c-471d-8ae0-9bd60d70567b.60b44875e55860e846c7f9475838e5c52eca4020.zh-cn.xlf")

# --- de-de sheet: Status column + new Target/Handback file links ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("H2").Value = "2016-03-17 16:11:42"
$dede.Range("H3").Value = "2016-03-17 16:11:42"

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/c2d708bfa23c3165cf7fa1ec8868c123901e405f/e2e/a1897a40-c0ff-4968-99fa-885e084b0050.md", "", "", "a1897a40-c0ff-4968-99fa-885e084b0050.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f042856f2252775f163621ae6847f0f6dbfb0c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/a1897a40-c0ff-4968-99fa-885e084b0050.0ffc24064d9e76e357c05896c12c1906504eaec3.de-de.xlf", "", "", "a1897a40-c0ff-4968-99fa-885e084b0050.0ffc24064d9e76e357c05896c12c1906504eaec3.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/c2d708bfa23c3165cf7fa1ec8868c123901e405f/e2e/e6dee80d-419c-471d-8ae0-9bd60d70567b.md", "", "", "e6dee80d-419c-471d-8ae0-9bd60d70567b.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f042856f2252775f163621ae6847f0f6dbfb0c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/e6dee80d-419c-471d-8ae0-9bd60d70567b.60b44875e55860e846c7f9475838e5c52eca4020.de-de.xlf", "", "", "e6dee80d-419c-471d-8ae0-9bd60d70567b.60b44875e55860e846c7f9475838e5c52eca4020.de-de.xlf")
